$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new "Region" column before the current Date column (col B),
# pushing Date/CFR2/DiffDE/AgeCompDE/RateCompDE/relAgeDE/relRateDE one
# column to the right (B->C, C->D, ... H->I). This also shifts the
# custom column-width formatting from col B to col C automatically.
$ws.Columns.Item(2).Insert()

# Add a new row at the bottom of the data (row 9) for the extra
# USA/NYC observation.
$ws.Rows.Item(9).Insert()

# Header row
$ws.Range("B1").Value = "Region"

$data = @(
    @("USA", "All", 43906, 0.01796725060231124, 0.01408418784366104, -0.004529216959880018, 0.01861340480354106, 0.1957088961735027, 0.8042911038264974),
    @("SouthKorea", "All", 43943, 0.02225547035720965, 0.009795968088762631, 0.01243792488709594, -0.002641956798333309, 0.8248025512769064, 0.1751974487230937),
    @("China", "All", 43872, 0.02290025071633238, 0.009151187729639902, 0.007221859156972936, 0.001929328572666963, 0.7891717851642323, 0.2108282148357678),
    @("Germany", "All", 43941, 0.03205143844597228, 0, 0, 0, $null, $null),
    @("France", "All", 43914, 0.03983587515221891, -0.007784436706246629, -0.001800858098157913, -0.005983578608088725, 0.2313408363527204, 0.7686591636472796),
    @("USA", "NYC", 43941, 0.0708957990420689, -0.03884436059609662, 0.00714009319979759, -0.04598445379589421, 0.1344029004214686, 0.8655970995785315),
    @("Spain", "All", 43937, 0.1050210003716739, -0.07296956192570159, -0.0331076802397085, -0.0398618816859931, 0.4537190489565924, 0.5462809510434077),
    @("Italy", "All", 43941, 0.1272752828730058, -0.09522384442703355, -0.04483297587587955, -0.050390868551154, 0.4708166966545168, 0.5291833033454833)
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    if ($row[7] -ne $null) {
        $ws.Cells.Item($r, 8).Value = $row[7]
    }
    if ($row[8] -ne $null) {
        $ws.Cells.Item($r, 9).Value = $row[8]
    }
    $r = $r + 1
}
